$wb = $excel.ActiveWorkbook

# --- Overview sheet: row 3 (b6121ff3 file) status columns B and C ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("B3").Value = "Handed back: in sync with en-US"
$wsOverview.Range("C3").Value = "Handed back: in sync with en-US"

# --- zh-cn sheet: row 2 (2aedfd63) Latest Handback DateTime, row 3 (b6121ff3) Status + Latest Handback DateTime ---
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("H2").Value = "2016-03-19 04:11:01"
$wsZhCn.Range("C3").Value = "Handed back: in sync with en-US"
$wsZhCn.Range("H3").Value = "2016-03-19 04:11:01"

# --- de-de sheet: row 2 (2aedfd63) Latest Handback DateTime, row 3 (b6121ff3) Status + Latest Handback DateTime ---
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("H2").Value = "2016-03-19 04:11:15"
$wsDeDe.Range("C3").Value = "Handed back: in sync with en-US"
$wsDeDe.Range("H3").Value = "2016-03-19 04:11:15"
